$d = $word.ActiveDocument

# Helper: append a bold " " run followed by a bold+green "[DONE]" run
# to the end of the given paragraph (before its end-of-paragraph mark).
function Add-DoneMarker($para) {
    $r = $para.Range
    $insStart = $r.End - 1
    $r.InsertAfter(" ")
    $spaceRange = $d.Range($insStart, $r.End - 1)
    $spaceRange.Font.Bold = $true

    $insStart2 = $r.End - 1
    $r.InsertAfter("[DONE]")
    $doneRange = $d.Range($insStart2, $r.End - 1)
    $doneRange.Font.Bold = $true
    $doneRange.Font.Color = 32768
}

# ---------------------------------------------------------------------
# 1) "Load conditions used for benchmarking ... /min?" paragraph:
#    append " [DONE]" (bold, green) at the end.
# ---------------------------------------------------------------------
$pLoad = $d.Paragraphs(15)
Add-DoneMarker $pLoad

# ---------------------------------------------------------------------
# 2) "For private cloud offerings ... applications." paragraph:
#    merge the three runs that make up the sentence into a single run,
#    leaving the trailing " [DONE]" runs untouched.
# ---------------------------------------------------------------------
$pPrivate = $d.Paragraphs(31)
$fullText = $pPrivate.Range.Text
$marker = "applications."
$appIdx = $fullText.IndexOf($marker)
$endOfSentence = $pPrivate.Range.Start + $appIdx + $marker.Length
$sentenceRange = $d.Range($pPrivate.Range.Start, $endOfSentence)
$insPoint = $sentenceRange.Start
$newSentence = "For private cloud offerings, where operators may have no limitation of how long a service can run, we may see different execution paths appearing in the applications."
$sentenceRange.Delete()
$insertionRange = $d.Range($insPoint, $insPoint)
$insertionRange.InsertBefore($newSentence)
$finalSentenceRange = $d.Range($insPoint, $insPoint + $newSentence.Length)
$finalSentenceRange.Font.Bold = $true

# ---------------------------------------------------------------------
# 3) "What happens if there is a lot of variance ..." paragraph:
#    append " [DONE]" (bold, green) at the end.
# ---------------------------------------------------------------------
$pVariance = $d.Paragraphs(39)
Add-DoneMarker $pVariance

# ---------------------------------------------------------------------
# 4) "How to decide the duration of history data?" paragraph:
#    append " [DONE]" (bold, green) at the end.
# ---------------------------------------------------------------------
$pHistory = $d.Paragraphs(43)
Add-DoneMarker $pHistory

Write-Output "Done."
